# Add a new "clave" (password) column to the parameters sheet, mirroring
# the existing header / data formatting already used in columns A-C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + value for column D
$ws.Range("D1").Value = "clave"
$ws.Range("D2").Value = "Scant9756"

# Match the formatting already used for the header row (bold, centered)
# and the plain (non-hyperlink) data row style.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

# Size column D to fit its new contents.
$ws.Columns.Item(4).ColumnWidth = 8.8

# Move/update the active selection, as left by the editor.
$ws.Range("C8").Select()
